# TileZoneTemplate.xlsx update
# Re-generated tile-zone coordinate data on Sheet1 (rows 2-31, columns A-D,
# plus one stray cell L3), the active-cell selection, and the bestFit
# column widths that widened to fit the refreshed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Updated data values (Cave_x/Cave_y/Mountains_x/Mountains_y columns, etc.) ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 12

$ws.Range("B3").Value = 2
$ws.Range("D3").Value = 10
$ws.Range("L3").Value = 11

$ws.Range("B4").Value = 3
$ws.Range("D4").Value = 12

$ws.Range("B5").Value = 4
$ws.Range("D5").Value = 13

$ws.Range("B6").Value = 5
$ws.Range("D6").Value = 14

$ws.Range("A7").Value = 2
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 15

$ws.Range("A8").Value = 3
$ws.Range("B8").Value = 6
$ws.Range("D8").Value = 9

$ws.Range("B9").Value = 2
$ws.Range("D9").Value = 10

$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 12

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 15

$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 16

$ws.Range("A13").Value = 4
$ws.Range("B13").Value = 6
$ws.Range("D13").Value = 10

$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 4
$ws.Range("D14").Value = 12

$ws.Range("A15").Value = 5
$ws.Range("B15").Value = 6
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 14

$ws.Range("B16").Value = 1
$ws.Range("C16").Value = 4

$ws.Range("B17").Value = 2
$ws.Range("C17").Value = 5

$ws.Range("A18").Value = 6
$ws.Range("B18").Value = 4
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 12

$ws.Range("A19").Value = 6
$ws.Range("B19").Value = 6
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 15

$ws.Range("B20").Value = 2
$ws.Range("D20").Value = 10

$ws.Range("A21").Value = 7
$ws.Range("B21").Value = 4
$ws.Range("D21").Value = 11

$ws.Range("A22").Value = 7
$ws.Range("B22").Value = 6
$ws.Range("D22").Value = 12

$ws.Range("B23").Value = 2
$ws.Range("C23").Value = 6

$ws.Range("B24").Value = 4
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 14

$ws.Range("A25").Value = 8
$ws.Range("B25").Value = 5
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 15

$ws.Range("A26").Value = 8
$ws.Range("B26").Value = 6
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 13

$ws.Range("D27").Value = 10

$ws.Range("D28").Value = 11

$ws.Range("C29").Value = 8

$ws.Range("C30").Value = 8
$ws.Range("D30").Value = 13

$ws.Range("C31").Value = 8
$ws.Range("D31").Value = 15

# --- Column widths widened (bestFit) to accommodate the refreshed data ---
$ws.Columns.Item(1).ColumnWidth = 9.333333333333334
$ws.Columns.Item(2).ColumnWidth = 9.333333333333334
$ws.Columns.Item(3).ColumnWidth = 15.166666666666666
$ws.Columns.Item(4).ColumnWidth = 15.166666666666666
$ws.Columns.Item(5).ColumnWidth = 11.333333333333334
$ws.Columns.Item(6).ColumnWidth = 11.333333333333334
$ws.Columns.Item(7).ColumnWidth = 19.666666666666668
$ws.Columns.Item(8).ColumnWidth = 19.666666666666668
$ws.Columns.Item(9).ColumnWidth = 26
$ws.Columns.Item(10).ColumnWidth = 26
$ws.Columns.Item(11).ColumnWidth = 22
$ws.Columns.Item(12).ColumnWidth = 22

# --- Active cell / selection moved to U31 ---
$ws.Range("U31").Select()
